$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B88 was stored as a text string "1"; change it to a real number 1
$ws.Range("B88").Value = 1

# Append new row 89 with the new annotation data
$ws.Range("A89").Value = "Ying Tang"
# B89 must stay a text "4" (not a number). A plain Value="4" assignment would
# be auto-converted to a number, so build it as a text formula in a scratch
# cell and paste-special just the value/type across, then clean up the
# scratch cell.
$ws.Range("Z1").Formula = '="4"'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("B89").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$ws.Range("Z1").Clear() | Out-Null
$excel.CutCopyMode = 0
$ws.Range("C89").Value = "will be helpful"
$ws.Range("D89").Value = "SUG"
$ws.Range("E89").Value = "MET"
$ws.Range("F89").Value = "835bb598-ab4d-49cf-8a20-53efcbaa448c"
$ws.Range("G89").Value = "rJSr0GZR-_annotated.xlsx"
$ws.Range("H89").Value = "Maybe some visualization like t-sne will be helpful."
